$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns AE:AJ on row 1 (copy formatting from AD1) ---
$ws.Range("AD1").Copy()
$ws.Range("AE1:AJ1").PasteSpecial(-4122)

$ws.Range("AE1").Value = "Kutsumanimi"
$ws.Range("AF1").Value = "Syntymäpaikka"
$ws.Range("AG1").Value = "Passin numero"
$ws.Range("AH1").Value = "Kansallinen ID-tunnus"
$ws.Range("AI1").Value = "Kaupunki ja maa"
$ws.Range("AJ1").Value = "Hakemus-oid"

# --- Update existing hakija OID sample value ---
$ws.Range("G2").Value = "Hakijaoid1"

# --- New data values for row 2 (copy formatting from AD2) ---
$ws.Range("AD2").Copy()
$ws.Range("AE2:AJ2").PasteSpecial(-4122)

$ws.Range("AE2").Value = "Tuomas"
$ws.Range("AJ2").Value = "Hakemus1"

# --- Extend the blank bordered block (rows 3:10) across the new columns ---
$ws.Range("AD3").Copy()
$ws.Range("AE3:AJ10").PasteSpecial(-4122)

$excel.CutCopyMode = 0
